$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "ValidLogin"

# Set data row first (manager needs to land right after admin in shared strings)
$ws.Range("B2").Value = "manager"

# Set header row (UserName/Password appended after)
$ws.Range("A1").Value = "UserName"
$ws.Range("B1").Value = "Password"

# Keep A2 as admin (already present)
$ws.Range("A2").Value = "admin"

# Best-fit the two used columns to their content (matches bestFit columns in target)
$ws.Columns("A:B").AutoFit()

# Set selection to C4 (as shown in diff)
$ws.Range("C4").Select()
